# Add upload form tests: fill in a sample form_title/form_id row on the
# "settings" sheet and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Make "settings" the active/tab-selected sheet (it was already the tab
# selected sheet in the source workbook).
$ws.Activate()

# Row 1 already holds the headers "form_title" / "form_id"; populate row 2
# with sample values used by the new upload form tests.
$ws.Range("A2").Value = "Just a test"
$ws.Range("B2").Value = "Justtest-date"

# Move/record the selection on the sheet as captured in the saved file.
$ws.Range("B3").Select()
